$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Metairie"
$ws.Range("C4").Value = "San Jose"
$ws.Range("C6").Value = "Portsmouth"
$ws.Range("C8").Value = "Oak Park"
$ws.Range("A10").Value = "Carol Burnett"
$ws.Range("C10").Value = "San Antonio"
$ws.Range("C13").Value = "West Germany"
$ws.Range("C16").Value = "Andrews"
$ws.Range("C19").Value = "Fritch"
$ws.Range("C21").Value = "Pawnee City"
$ws.Range("C23").Value = "Brooklyn"
$ws.Range("C25").Value = "Denver"
$ws.Range("C27").Value = "Missoula"
